$d = $word.ActiveDocument

# --- 1. Merge split text runs for Title, Author and Abstract paragraphs ---
# The content does not change, but the runs that were previously split on
# every word/space get consolidated into a single <w:r> per paragraph, exactly
# as in the target revision. We replace each paragraph (via InsertXML on its
# full Range, including the paragraph mark) with an equivalent paragraph that
# has a single run, keeping the original pStyle and xml:space="preserve".

$titleXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Questions: Integration by substitution</w:t></w:r></w:p>
'@
$d.Paragraphs(1).Range.InsertXML($titleXml) | Out-Null

$authorXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">Donald Campbell</w:t></w:r></w:p>
'@
$d.Paragraphs(2).Range.InsertXML($authorXml) | Out-Null

$abstractXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Abstract"/></w:pPr><w:r><w:t xml:space="preserve">A selection of questions for the study guide on integration by substitution.</w:t></w:r></w:p>
'@
$d.Paragraphs(4).Range.InsertXML($abstractXml) | Out-Null

# --- 2. Reorder <m:sepChr> before <m:endChr> inside every <m:dPr> of every bracket ---
#        delimiter used across the equations (m:d). Word stores the delimiter
#        characters as begChr / sepChr / endChr and the canonical order produced
#        by a fresh edit places sepChr immediately after begChr.

$omathFixes = @(
    @{ Index = 2; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 3; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>5</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 4; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>x</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>4</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 5; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 6; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>4</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 7; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>7</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 8; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>x</m:t></m:r></m:num><m:den><m:r><m:t>5</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 9; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>1</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 10; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 11; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>6</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 13; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 14; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 15; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>5</m:t></m:r></m:num><m:den><m:r><m:t>6</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 16; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 17; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>x</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 18; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>4</m:t></m:r></m:num><m:den><m:r><m:t>5</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>π</m:t></m:r></m:num><m:den><m:r><m:t>4</m:t></m:r></m:den></m:f></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 19; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>π</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>4</m:t></m:r><m:r><m:t>x</m:t></m:r></m:num><m:den><m:r><m:t>9</m:t></m:r></m:den></m:f></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 20; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>π</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 21; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:t>4</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>x</m:t></m:r></m:num><m:den><m:r><m:t>4</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>π</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 22; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>5</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>π</m:t></m:r></m:num><m:den><m:r><m:t>6</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 26; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:d></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 27; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>exp</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>x</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 31; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>5</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>1</m:t></m:r></m:den></m:f><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 34; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:t>6</m:t></m:r><m:r><m:t>x</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>4</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 35; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:t>5</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>7</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 36; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:t>8</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>exp</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 37; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>1</m:t></m:r></m:num><m:den><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:den></m:f><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 38; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:t>6</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 39; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>exp</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 40; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>x</m:t></m:r></m:num><m:den><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>3</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:den></m:f><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 42; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>4</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ },
    @{ Index = 43; Xml = @'
<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:num><m:den><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:den></m:f><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath>
'@ }
)

foreach ($fix in $omathFixes) {
    $om = $d.OMaths.Item($fix.Index)
    $om.Range.InsertXML($fix.Xml) | Out-Null
}

Write-Host "Done. Paragraphs:" $d.Paragraphs.Count "OMaths:" $d.OMaths.Count